$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.6151236666666667
$ws.Cells.Item(2, 8).Value = 1.845371
$ws.Cells.Item(2, 9).Value = 0.01505763605988265
$ws.Cells.Item(2, 10).Value = 0.01505763605988265
$ws.Cells.Item(2, 13).Value = 3.390903999999999
$ws.Cells.Item(2, 14).Value = 10.172712
$ws.Cells.Item(2, 15).Value = 0.1656941395696903
$ws.Cells.Item(2, 16).Value = 0.1656941395696903
$ws.Cells.Item(2, 17).Value = 2.085825301794666
$ws.Cells.Item(2, 18).Value = 18.772427716152
$ws.Cells.Item(2, 19).Value = 0.002494962050895797
$ws.Cells.Item(2, 20).Value = 0.002494962050895797

$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.6151236666666667
$ws.Cells.Item(3, 8).Value = 1.845371
$ws.Cells.Item(3, 9).Value = 0.01505763605988265
$ws.Cells.Item(3, 10).Value = 0.01505763605988265
$ws.Cells.Item(3, 15).Value = 0.3503621390046235
$ws.Cells.Item(3, 16).Value = 0.3503621390046235
$ws.Cells.Item(3, 17).Value = 4.410501277985
$ws.Cells.Item(3, 18).Value = 39.694511501865
$ws.Cells.Item(3, 19).Value = 0.005275625578293638
$ws.Cells.Item(3, 20).Value = 0.005275625578293637

$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.6151236666666667
$ws.Cells.Item(4, 8).Value = 1.845371
$ws.Cells.Item(4, 9).Value = 0.01505763605988265
$ws.Cells.Item(4, 10).Value = 0.01505763605988265
$ws.Cells.Item(4, 13).Value = 0.9786293333333332
$ws.Cells.Item(4, 14).Value = 2.935888
$ws.Cells.Item(4, 15).Value = 0.04782003422813688
$ws.Cells.Item(4, 16).Value = 0.04782003422813689
$ws.Cells.Item(4, 17).Value = 0.6019780638275556
$ws.Cells.Item(4, 18).Value = 5.417802574448
$ws.Cells.Item(4, 19).Value = 0.0007200566717784166
$ws.Cells.Item(4, 20).Value = 0.0007200566717784166

$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.6151236666666667
$ws.Cells.Item(5, 8).Value = 1.845371
$ws.Cells.Item(5, 9).Value = 0.01505763605988265
$ws.Cells.Item(5, 10).Value = 0.01505763605988265
$ws.Cells.Item(5, 13).Value = 5.160571999999999
$ws.Cells.Item(5, 14).Value = 15.481716
$ws.Cells.Item(5, 15).Value = 0.2521677220078881
$ws.Cells.Item(5, 16).Value = 0.2521677220078881
$ws.Cells.Item(5, 17).Value = 3.174389970737333
$ws.Cells.Item(5, 18).Value = 28.569509736636
$ws.Cells.Item(5, 19).Value = 0.003797049784044439
$ws.Cells.Item(5, 20).Value = 0.003797049784044439

$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.6151236666666667
$ws.Cells.Item(6, 8).Value = 1.845371
$ws.Cells.Item(6, 9).Value = 0.01505763605988265
$ws.Cells.Item(6, 10).Value = 0.01505763605988265
$ws.Cells.Item(6, 13).Value = 3.764629333333333
$ws.Cells.Item(6, 14).Value = 11.293888
$ws.Cells.Item(6, 15).Value = 0.1839559651896613
$ws.Cells.Item(6, 16).Value = 0.1839559651896613
$ws.Cells.Item(6, 17).Value = 2.315712599160889
$ws.Cells.Item(6, 18).Value = 20.841413392448
$ws.Cells.Item(6, 19).Value = 0.002769941974870363
$ws.Cells.Item(6, 20).Value = 0.002769941974870362

$ws.Cells.Item(7, 9).Value = 0.9124571722898065
$ws.Cells.Item(7, 10).Value = 0.9124571722898065
$ws.Cells.Item(7, 13).Value = 3.390903999999999
$ws.Cells.Item(7, 14).Value = 10.172712
$ws.Cells.Item(7, 15).Value = 0.1656941395696903
$ws.Cells.Item(7, 16).Value = 0.1656941395696903
$ws.Cells.Item(7, 17).Value = 126.396085627064
$ws.Cells.Item(7, 18).Value = 1137.564770643576
$ws.Cells.Item(7, 19).Value = 0.1511888060567521
$ws.Cells.Item(7, 20).Value = 0.1511888060567521

$ws.Cells.Item(8, 9).Value = 0.9124571722898065
$ws.Cells.Item(8, 10).Value = 0.9124571722898065
$ws.Cells.Item(8, 15).Value = 0.3503621390046235
$ws.Cells.Item(8, 16).Value = 0.3503621390046235
$ws.Cells.Item(8, 19).Value = 0.3196904466335669
$ws.Cells.Item(8, 20).Value = 0.3196904466335669

$ws.Cells.Item(9, 9).Value = 0.9124571722898065
$ws.Cells.Item(9, 10).Value = 0.9124571722898065
$ws.Cells.Item(9, 13).Value = 0.9786293333333332
$ws.Cells.Item(9, 14).Value = 2.935888
$ws.Cells.Item(9, 15).Value = 0.04782003422813688
$ws.Cells.Item(9, 16).Value = 0.04782003422813689
$ws.Cells.Item(9, 17).Value = 36.47844852380265
$ws.Cells.Item(9, 18).Value = 328.306036714224
$ws.Cells.Item(9, 19).Value = 0.04363373321060754
$ws.Cells.Item(9, 20).Value = 0.04363373321060755

$ws.Cells.Item(10, 9).Value = 0.9124571722898065
$ws.Cells.Item(10, 10).Value = 0.9124571722898065
$ws.Cells.Item(10, 13).Value = 5.160571999999999
$ws.Cells.Item(10, 14).Value = 15.481716
$ws.Cells.Item(10, 15).Value = 0.2521677220078881
$ws.Cells.Item(10, 16).Value = 0.2521677220078881
$ws.Cells.Item(10, 17).Value = 192.360532883452
$ws.Cells.Item(10, 18).Value = 1731.244795951068
$ws.Cells.Item(10, 19).Value = 0.2300922465660795
$ws.Cells.Item(10, 20).Value = 0.2300922465660795

$ws.Cells.Item(11, 9).Value = 0.9124571722898065
$ws.Cells.Item(11, 10).Value = 0.9124571722898065
$ws.Cells.Item(11, 13).Value = 3.764629333333333
$ws.Cells.Item(11, 14).Value = 11.293888
$ws.Cells.Item(11, 15).Value = 0.1839559651896613
$ws.Cells.Item(11, 16).Value = 0.1839559651896613
$ws.Cells.Item(11, 17).Value = 140.3267127498027
$ws.Cells.Item(11, 18).Value = 1262.940414748224
$ws.Cells.Item(11, 19).Value = 0.1678519398228005
$ws.Cells.Item(11, 20).Value = 0.1678519398228005

$ws.Cells.Item(12, 7).Value = 2.961112666666666
$ws.Cells.Item(12, 8).Value = 8.883337999999998
$ws.Cells.Item(12, 9).Value = 0.07248519165031087
$ws.Cells.Item(12, 10).Value = 0.07248519165031085
$ws.Cells.Item(12, 13).Value = 3.390903999999999
$ws.Cells.Item(12, 14).Value = 10.172712
$ws.Cells.Item(12, 15).Value = 0.1656941395696903
$ws.Cells.Item(12, 16).Value = 0.1656941395696903
$ws.Cells.Item(12, 17).Value = 10.04084878585066
$ws.Cells.Item(12, 18).Value = 90.36763907265598
$ws.Cells.Item(12, 19).Value = 0.01201037146204235
$ws.Cells.Item(12, 20).Value = 0.01201037146204235

$ws.Cells.Item(13, 7).Value = 2.961112666666666
$ws.Cells.Item(13, 8).Value = 8.883337999999998
$ws.Cells.Item(13, 9).Value = 0.07248519165031087
$ws.Cells.Item(13, 10).Value = 0.07248519165031085
$ws.Cells.Item(13, 15).Value = 0.3503621390046235
$ws.Cells.Item(13, 16).Value = 0.3503621390046235
$ws.Cells.Item(13, 17).Value = 21.23148873682999
$ws.Cells.Item(13, 18).Value = 191.0833986314699
$ws.Cells.Item(13, 19).Value = 0.02539606679276299
$ws.Cells.Item(13, 20).Value = 0.02539606679276299

$ws.Cells.Item(14, 7).Value = 2.961112666666666
$ws.Cells.Item(14, 8).Value = 8.883337999999998
$ws.Cells.Item(14, 9).Value = 0.07248519165031087
$ws.Cells.Item(14, 10).Value = 0.07248519165031085
$ws.Cells.Item(14, 13).Value = 0.9786293333333332
$ws.Cells.Item(14, 14).Value = 2.935888
$ws.Cells.Item(14, 15).Value = 0.04782003422813688
$ws.Cells.Item(14, 16).Value = 0.04782003422813689
$ws.Cells.Item(14, 17).Value = 2.897831714904888
$ws.Cells.Item(14, 18).Value = 26.08048543414399
$ws.Cells.Item(14, 19).Value = 0.003466244345750927
$ws.Cells.Item(14, 20).Value = 0.003466244345750927

$ws.Cells.Item(15, 7).Value = 2.961112666666666
$ws.Cells.Item(15, 8).Value = 8.883337999999998
$ws.Cells.Item(15, 9).Value = 0.07248519165031087
$ws.Cells.Item(15, 10).Value = 0.07248519165031085
$ws.Cells.Item(15, 13).Value = 5.160571999999999
$ws.Cells.Item(15, 14).Value = 15.481716
$ws.Cells.Item(15, 15).Value = 0.2521677220078881
$ws.Cells.Item(15, 16).Value = 0.2521677220078881
$ws.Cells.Item(15, 17).Value = 15.28103511644533
$ws.Cells.Item(15, 18).Value = 137.529316048008
$ws.Cells.Item(15, 19).Value = 0.01827842565776408
$ws.Cells.Item(15, 20).Value = 0.01827842565776408

$ws.Cells.Item(16, 7).Value = 2.961112666666666
$ws.Cells.Item(16, 8).Value = 8.883337999999998
$ws.Cells.Item(16, 9).Value = 0.07248519165031087
$ws.Cells.Item(16, 10).Value = 0.07248519165031085
$ws.Cells.Item(16, 13).Value = 3.764629333333333
$ws.Cells.Item(16, 14).Value = 11.293888
$ws.Cells.Item(16, 15).Value = 0.1839559651896613
$ws.Cells.Item(16, 16).Value = 0.1839559651896613
$ws.Cells.Item(16, 17).Value = 11.14749160423822
$ws.Cells.Item(16, 18).Value = 100.327424438144
$ws.Cells.Item(16, 19).Value = 0.01333408339199052
$ws.Cells.Item(16, 20).Value = 0.01333408339199051
